$wb = $excel.ActiveWorkbook
$wsVeda = $wb.Worksheets.Item("Veda")
$wsHist = $wb.Worksheets.Item("historical_data")

# Update the "input cell" marker label to reflect the single changed input (B3)
$wsVeda.Range("B2").Value = "~Inputcell: 3"

# Change the scenario selector from 2 (Net Zero 2050) to 3 (Nationally Determined Contributions (NDCs))
$wsVeda.Range("B3").Value = 3

# Reflect the user's on-screen navigation: selection moved on Veda, and the
# historical_data sheet became the active/selected tab.
$wsVeda.Range("Y26").Select()
$wsHist.Activate()
